# Applies "hybrid bold + color" highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) inside specific bullet/impact
# paragraphs, matching the commit's DOCX highlighting behavior: each metric
# becomes its own run with Bold + color 2C3E50, while the surrounding text
# stays in plain (unformatted) runs.

$HighlightColor = 5258796   # packed BGR int for RGB(0x2C,0x3E,0x50) -> w:color val="2C3E50"

# Locates the 1-based Paragraphs index of the paragraph whose text contains
# $fingerprint (a short, distinguishing substring of the target paragraph)
# but does NOT contain $exclude (used to disambiguate two paragraphs where
# one's text is a prefix of the other's, e.g. the short vs. long "Achieved
# 87%..." bullets).
function Find-ParagraphIndex($doc, [string]$fingerprint, [string]$exclude) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        if ($t.Contains($fingerprint)) {
            if ([string]::IsNullOrEmpty($exclude) -or -not $t.Contains($exclude)) {
                return $i
            }
        }
    }
    return -1
}

# Within the paragraph uniquely identified by $fingerprint (optionally minus
# $exclude matches), bolds + colors each substring listed in $targets
# (matched in order, left to right, so repeated values like two separate
# "71%" runs are each found at their correct, advancing position).
function Apply-MetricHighlight($doc, [string]$fingerprint, [string[]]$targets, [string]$exclude = '') {
    $paraIndex = Find-ParagraphIndex $doc $fingerprint $exclude
    if ($paraIndex -lt 0) {
        Write-Output "PARAGRAPH NOT FOUND for fingerprint: $fingerprint"
        return
    }
    $p = $doc.Paragraphs.Item($paraIndex)
    $pStart = $p.Range.Start
    $full = $p.Range.Text
    $cursor = 0
    foreach ($t in $targets) {
        $idx = $full.IndexOf($t, $cursor)
        if ($idx -lt 0) {
            Write-Output "NOT FOUND: '$t' in paragraph index $paraIndex"
            continue
        }
        $rStart = $pStart + $idx
        $rEnd = $rStart + $t.Length
        $r = $doc.Range($rStart, $rEnd)
        $r.Font.Bold = 1
        $r.Font.Color = $HighlightColor
        $cursor = $idx + $t.Length
    }
}

$d = $word.ActiveDocument

# "• Discovered systematic race coding errors affecting all Black and Asian-American voters,
#   developed geospatial machine learning algorithms improving demographic classification
#   accuracy from 23% to 64%"
Apply-MetricHighlight $d 'Discovered systematic race coding errors' @('23%', '64%')

# "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%,
#   reducing polling error margins from ±4.2% to ±2.1%"
Apply-MetricHighlight $d 'reducing polling error margins' @('87%', '71%', '±4.2%', '±2.1%')

# "• Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
Apply-MetricHighlight $d 'Wrote RFP and analyzed bids' @('1,200')

# "• Created comprehensive meta-analysis framework handling millions of survey responses
#   that became the $400M Polling Consortium Database at The Analyst Institute, now valued
#   at $1B+"
Apply-MetricHighlight $d 'Polling Consortium Database' @('$400M', '$1B')

# "• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
Apply-MetricHighlight $d 'Algorithm reduced mapping costs' @('73.5%', '$4.7M')

# "• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
# (short-form achievement bullet; exclude the long-form one above which
# shares this exact text as a prefix)
Apply-MetricHighlight $d 'Achieved 87% prediction accuracy' @('87%', '71%') 'reducing polling error margins'

Write-Output "DONE"
